# Apply "added Tree Module with Try editor and excel data for dataprovider"
#
# - fills in the (previously empty) "tryeditor" sheet with a small code/error
#   table
# - appends a new data row to the "register" sheet
# - leaves "tryeditor" as the active/selected tab instead of "login"

$wb = $excel.ActiveWorkbook

$wsRegister  = $wb.Worksheets.Item(4)   # "register"
$wsTryEditor = $wb.Worksheets.Item(5)   # "tryeditor"

# --- tryeditor sheet data -------------------------------------------------
# Cells are written column-by-column (not row-by-row) so that the shared
# string table is populated in the same order it originally was authored.

$wsTryEditor.Range("A1").Value = "code"
$wsTryEditor.Range("A2").Value = "print(`"Hello"
$wsTryEditor.Range("A3").Value = "edgr7&rk"

$wsTryEditor.Range("B2").Value = "SyntaxError: bad input on line 1"
$wsTryEditor.Range("B3").Value = "NameError: name 'edgr7' is not defined on line 1"
$wsTryEditor.Range("B1").Value = "alert "

$wsTryEditor.Range("C2").Value = "syntax"
$wsTryEditor.Range("C3").Value = "name"
$wsTryEditor.Range("C1").Value = "error"

# Column widths (approximate autofit sizing for the new columns)
$wsTryEditor.Columns.Item(1).ColumnWidth = 21.666666666666668
$wsTryEditor.Columns.Item(2).ColumnWidth = 43.833333333333336
$wsTryEditor.Columns.Item(3).ColumnWidth = 14.666666666666666

# --- register sheet: new data row -----------------------------------------
$wsRegister.Range("A11").Value = "edgr7&rk"
$wsRegister.Range("B11").Value = "NameError: name 'edgr7' is not defined on line 1"

# --- selection / active tab -------------------------------------------------
# Selecting on "register" then "tryeditor" last leaves "tryeditor" as the
# active sheet/tab (clearing the previous "login" tabSelected flag).
$wsRegister.Range("B11").Select()
$wsTryEditor.Range("C1").Select()
